$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix F2: rank was stored as the shared string "Hanshi, 10"; should be numeric 10 ---
$ws.Range("F2").Value = 10

# --- Append new roster entries (rows 30-48) ---
# Row 30
$ws.Range("A30").Value = 'Jim'
$ws.Range("B30").Value = 'Kass'
$ws.Range("C30").Value = 'Alaska Shido-kan'
$ws.Range("D30").Value = 'Kenai, Alaska'
$ws.Range("E30").Value = 'USA'
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 'dan'
$ws.Range("H30").Value = 'jimkass.JPG'

# Row 31
$ws.Range("A31").Value = 'Sandy'
$ws.Range("B31").Value = 'Kass'
$ws.Range("C31").Value = 'Alaska Shido-kan'
$ws.Range("D31").Value = 'Kenai, Alaska'
$ws.Range("E31").Value = 'USA'
$ws.Range("H31").Value = 'sandykass.jpg'

# Row 32
$ws.Range("A32").Value = 'Keeven'
$ws.Range("B32").Value = 'Macik'
$ws.Range("C32").Value = 'Alaska Shido-kan'
$ws.Range("D32").Value = 'Kenai, Alaska'
$ws.Range("E32").Value = 'USA'
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 'dan'
$ws.Range("H32").Value = 'keevenmacik.JPG'

# Row 33
$ws.Range("A33").Value = 'Hannele'
$ws.Range("B33").Value = 'Zubeck'
$ws.Range("C33").Value = 'Alaska Shido-kan'
$ws.Range("D33").Value = 'Kenai, Alaska'
$ws.Range("E33").Value = 'USA'
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 'dan'
$ws.Range("H33").Value = 'hannelezubeck.jpg'

# Row 34
$ws.Range("A34").Value = 'Maija'
$ws.Range("B34").Value = 'Zubeck'
$ws.Range("C34").Value = 'Alaska Shido-kan'
$ws.Range("D34").Value = 'Kenai, Alaska'
$ws.Range("E34").Value = 'USA'
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 'dan'
$ws.Range("H34").Value = 'maijazubeck.JPG'

# Row 35
$ws.Range("A35").Value = 'Roberto'
$ws.Range("B35").Value = 'Curtis'
$ws.Range("C35").Value = 'Virginia'
$ws.Range("D35").Value = 'Alexandria, Virginia'
$ws.Range("E35").Value = 'USA'
$ws.Range("F35").Value = 8
$ws.Range("G35").Value = 'dan'

# Row 36
$ws.Range("A36").Value = 'Robert'
$ws.Range("B36").Value = 'Pangelinan'
$ws.Range("C36").Value = 'Virginia'
$ws.Range("D36").Value = 'Lorton, Virginia'
$ws.Range("E36").Value = 'USA'
$ws.Range("F36").Value = 6
$ws.Range("G36").Value = 'dan'

# Row 37
$ws.Range("A37").Value = 'Maria'
$ws.Range("B37").Value = 'Pangelinan'
$ws.Range("C37").Value = 'Virginia'
$ws.Range("D37").Value = 'Lorton, Virginia'
$ws.Range("E37").Value = 'USA'
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 'dan'

# Row 38
$ws.Range("A38").Value = 'Arthur'
$ws.Range("B38").Value = 'Flax'
$ws.Range("C38").Value = 'Virginia'
$ws.Range("D38").Value = 'Gaithersburg, Maryland'
$ws.Range("E38").Value = 'USA'
$ws.Range("F38").Value = 3
$ws.Range("G38").Value = 'dan'

# Row 39
$ws.Range("A39").Value = 'Denise'
$ws.Range("B39").Value = 'Coursey'
$ws.Range("C39").Value = 'Virginia'
$ws.Range("D39").Value = 'Alexandria, Virginia'
$ws.Range("E39").Value = 'USA'
$ws.Range("F39").Value = 2
$ws.Range("G39").Value = 'dan'

# Row 40
$ws.Range("A40").Value = 'Matt'
$ws.Range("B40").Value = 'Coursey'
$ws.Range("C40").Value = 'Virginia'
$ws.Range("D40").Value = 'Alexandria, Virginia'
$ws.Range("E40").Value = 'USA'
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 'dan'

# Row 41
$ws.Range("A41").Value = 'Margot'
$ws.Range("B41").Value = 'Paz'
$ws.Range("C41").Value = 'Virginia'
$ws.Range("D41").Value = 'Alexandria, Virginia'
$ws.Range("E41").Value = 'USA'
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 'kyu'

# Row 42
$ws.Range("A42").Value = 'Mom'
$ws.Range("B42").Value = 'Paz'
$ws.Range("C42").Value = 'Virginia'
$ws.Range("D42").Value = 'Alexandria, Virginia'
$ws.Range("E42").Value = 'USA'

# Row 43
$ws.Range("A43").Value = 'Noell'
$ws.Range("B43").Value = 'Dunlap McMichael'
$ws.Range("C43").Value = 'Virginia'
$ws.Range("D43").Value = 'Arlington, Virginia'
$ws.Range("E43").Value = 'USA'
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 'kyu'

# Row 44
$ws.Range("A44").Value = 'Noah'
$ws.Range("B44").Value = 'McMichael'
$ws.Range("C44").Value = 'Virginia'
$ws.Range("D44").Value = 'Arlington, Virginia'
$ws.Range("E44").Value = 'USA'

# Row 45
$ws.Range("A45").Value = 'Daniel'
$ws.Range("B45").Value = 'Reese'
$ws.Range("C45").Value = 'Virginia'
$ws.Range("D45").Value = 'Burke, Virginia'
$ws.Range("E45").Value = 'USA'
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 'kyu'

# Row 46
$ws.Range("A46").Value = 'Daniel'
$ws.Range("B46").Value = 'Kumar'
$ws.Range("C46").Value = 'Virginia'
$ws.Range("D46").Value = 'Alexandria, Virginia'
$ws.Range("E46").Value = 'USA'
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 'kyu'

# Row 47
$ws.Range("A47").Value = 'Bryan'
$ws.Range("B47").Value = 'Hudson'
$ws.Range("C47").Value = 'Virginia'
$ws.Range("D47").Value = 'Alexandria, Virginia'
$ws.Range("E47").Value = 'USA'
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 'kyu'

# Row 48
$ws.Range("A48").Value = 'Sean'
$ws.Range("B48").Value = 'McCloskey'
$ws.Range("C48").Value = 'Virginia'
$ws.Range("D48").Value = 'Alexandria, Virginia'
$ws.Range("E48").Value = 'USA'
$ws.Range("F48").Value = 3
$ws.Range("G48").Value = 'kyu'

# --- Column width tweaks (approximate; engine rounds to its own pixel grid) ---
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(8).ColumnWidth = 24.333333333333332

# --- AutoFilter over the full data range ---
$ws.Range("A1:H48").AutoFilter()

# --- Hidden workbook-scoped defined name Excel writes alongside an AutoFilter ---
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=roster!`$A`$1:`$H`$48")
$filterName.Visible = $false

# --- Selection / view state ---
$ws.Range("D44").Select()